$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("CDCF-PMpPDOU")
$ws2.Range("B2").Formula = "=1/1.60934*10^12"
$ws2.Range("B10").Select()

$ws3 = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ws3.Range("B2").Formula = "=1/1.60934*10^12"
$ws3.Range("B25").Select()

$ws1 = $wb.Worksheets.Item("About")
$ws1.Select()
